$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("efa")

# --- Bold the "Raw" sub-header (row 4) and the trait-pair labels (rows 5-11) in column D ---
# D4 ("Raw") already has a bottom border + center alignment; just bold it.
$ws.Range("D4").Font.Bold = $true

# D5:D10 (trait-pair labels / the "-" formula cell) already centered, no border; bold them.
$ws.Range("D5:D10").Font.Bold = $true

# D11 (65.89, bottom border + 0.00 number format) gets bold too.
$ws.Range("D11").Font.Bold = $true

# --- New "Factor / Variance / Proportion / Cumulative" table in L5:O10 ---

# Header row (top+bottom thin border, like the existing n/With-Grit headers)
$ws.Range("L5").Value = "Factor"
$ws.Range("M5").Value = "Variance"
$ws.Range("N5").Value = "Proportion"
$ws.Range("O5").Value = "Cumulative"
$ws.Range("L5:O5").Borders.Item(8).Color = 0
$ws.Range("L5:O5").Borders.Item(8).LineStyle = 1
$ws.Range("L5:O5").Borders.Item(9).Color = 0
$ws.Range("L5:O5").Borders.Item(9).LineStyle = 1

# Factor 1..4 rows (no border)
$ws.Range("L6").Value = "Factor 1"
$ws.Range("M6").Value = 7.21736
$ws.Range("N6").Value = 0.2062
$ws.Range("O6").Formula = "=N6"

$ws.Range("L7").Value = "Factor 2"
$ws.Range("M7").Value = 6.33064
$ws.Range("N7").Value = 0.1809
$ws.Range("O7").Formula = "=O6+N7"

$ws.Range("L8").Value = "Factor 3"
$ws.Range("M8").Value = 4.12549
$ws.Range("N8").Value = 0.1179
$ws.Range("O8").Formula = "=O7+N8"

$ws.Range("L9").Value = "Factor 4"
$ws.Range("M9").Value = 2.99571
$ws.Range("N9").Value = 0.0856
$ws.Range("O9").Formula = "=O8+N9"

# Factor 5 row (bottom border only)
$ws.Range("L10").Value = "Factor 5"
$ws.Range("M10").Value = 2.39139
$ws.Range("N10").Value = 0.0683
$ws.Range("O10").Formula = "=O9+N10"
$ws.Range("L10:O10").Borders.Item(9).Color = 0
$ws.Range("L10:O10").Borders.Item(9).LineStyle = 1

# Number formatting for the Variance/Proportion/Cumulative columns
$ws.Range("M6:O10").NumberFormat = "0.00"

# Approximate "best fit" column widths for the new columns
$ws.Columns.Item(12).ColumnWidth = 7.02
$ws.Columns.Item(13).ColumnWidth = 7.88
$ws.Columns.Item(14).ColumnWidth = 9.74
$ws.Columns.Item(15).ColumnWidth = 10.31
